$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Unit Processes")

# Add a new "Coke Oven" unit process as row 9, matching the text style (s="1")
# used by the rest of the table (NumberFormat "@" / Text).
$newRow = $ws.Range("A9:H9")
$newRow.NumberFormat = "@"

# Values are assigned in this particular order so that new entries land
# in the shared-strings table in the same order as the target workbook
# (IEAGHG_coke_oven, coke oven, coke, SteelUnits_Variables.xlsx,
#  SteelUnits_Relationships.xlsx, Coke Oven).
$ws.Range("A9").Value = "IEAGHG_coke_oven"
$ws.Range("B9").Value = "coke oven"
$ws.Range("C9").Value = "coke"
$ws.Range("D9").Value = "output"
$ws.Range("E9").Value = "data/steel/SteelUnits_Variables.xlsx"
$ws.Range("G9").Value = "data/steel/SteelUnits_Relationships.xlsx"
$ws.Range("F9").Value = "Coke Oven"
$ws.Range("H9").Value = "Coke Oven"

$ws.Range("A2").Select()
